$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBE = New-Object 'object[,]' 24,4
$dataBE[0,0] = 0.2881169905109251
$dataBE[0,1] = 1.626987699542094
$dataBE[0,2] = 0.1496068669990043
$dataBE[0,3] = 0.5333859586016987
$dataBE[1,0] = 1.445647641019636
$dataBE[1,1] = 1.626987699542094
$dataBE[1,2] = 0.1496068669990043
$dataBE[1,3] = 0.5333859586016987
$dataBE[2,0] = 3.272327238179451
$dataBE[2,1] = 1.626987699542094
$dataBE[2,2] = 3.223369029078222
$dataBE[2,3] = 0.5333859586016987
$dataBE[3,0] = 3.272327238179451
$dataBE[3,1] = 1.626987699542094
$dataBE[3,2] = 3.223369029078222
$dataBE[3,3] = 0.5333859586016987
$dataBE[4,0] = 1.445647641019636
$dataBE[4,1] = 1.626987699542094
$dataBE[4,2] = 0.7210945179870265
$dataBE[4,3] = 0.5333859586016987
$dataBE[5,0] = 3.272327238179451
$dataBE[5,1] = 1.626987699542094
$dataBE[5,2] = 3.223369029078222
$dataBE[5,3] = 0.5333859586016987
$dataBE[6,0] = 3.272327238179451
$dataBE[6,1] = 1.626987699542094
$dataBE[6,2] = 0.7210945179870265
$dataBE[6,3] = 0.5333859586016987
$dataBE[7,0] = 0.1169995834814548
$dataBE[7,1] = 0.3048912486333797
$dataBE[7,2] = 0.7210945179870265
$dataBE[7,3] = 0.5333859586016987
$dataBE[8,0] = 0.6545652718822623
$dataBE[8,1] = 1.626987699542094
$dataBE[8,2] = 3.223369029078222
$dataBE[8,3] = 0.5333859586016987
$dataBE[9,0] = 3.272327238179451
$dataBE[9,1] = 1.626987699542094
$dataBE[9,2] = 3.223369029078222
$dataBE[9,3] = 0.5333859586016987
$dataBE[10,0] = 3.272327238179451
$dataBE[10,1] = 1.626987699542094
$dataBE[10,2] = 3.223369029078222
$dataBE[10,3] = 0.5333859586016987
$dataBE[11,0] = 3.272327238179451
$dataBE[11,1] = 0.3048912486333797
$dataBE[11,2] = 0.7210945179870265
$dataBE[11,3] = 0.5333859586016987
$dataBE[12,0] = 0.04172184405617529
$dataBE[12,1] = 0.3048912486333797
$dataBE[12,2] = 3.223369029078222
$dataBE[12,3] = 13.86384647080068
$dataBE[13,0] = 3.272327238179451
$dataBE[13,1] = 1.626987699542094
$dataBE[13,2] = 0.7210945179870265
$dataBE[13,3] = 0.5333859586016987
$dataBE[14,0] = 3.272327238179451
$dataBE[14,1] = 1.626987699542094
$dataBE[14,2] = 0.1496068669990043
$dataBE[14,3] = 0.5333859586016987
$dataBE[15,0] = 0.6545652718822623
$dataBE[15,1] = 1.626987699542094
$dataBE[15,2] = 0.7210945179870265
$dataBE[15,3] = 13.86384647080068
$dataBE[16,0] = 3.272327238179451
$dataBE[16,1] = 1.626987699542094
$dataBE[16,2] = 0.7210945179870265
$dataBE[16,3] = 0.5333859586016987
$dataBE[17,0] = 0.1169995834814548
$dataBE[17,1] = 1.626987699542094
$dataBE[17,2] = 0.7210945179870265
$dataBE[17,3] = 0.5333859586016987
$dataBE[18,0] = 3.272327238179451
$dataBE[18,1] = 1.626987699542094
$dataBE[18,2] = 0.1496068669990043
$dataBE[18,3] = 0.5333859586016987
$dataBE[19,0] = 1.445647641019636
$dataBE[19,1] = 1.626987699542094
$dataBE[19,2] = 0.7210945179870265
$dataBE[19,3] = 0.5333859586016987
$dataBE[20,0] = 3.272327238179451
$dataBE[20,1] = 1.626987699542094
$dataBE[20,2] = 0.1496068669990043
$dataBE[20,3] = 13.86384647080068
$dataBE[21,0] = 3.272327238179451
$dataBE[21,1] = 1.626987699542094
$dataBE[21,2] = 0.1496068669990043
$dataBE[21,3] = 13.86384647080068
$dataBE[22,0] = 3.272327238179451
$dataBE[22,1] = 1.626987699542094
$dataBE[22,2] = 0.7210945179870265
$dataBE[22,3] = 0.5333859586016987
$dataBE[23,0] = 3.272327238179451
$dataBE[23,1] = 1.626987699542094
$dataBE[23,2] = 3.223369029078222
$dataBE[23,3] = 0.5333859586016987

$dataG = New-Object 'object[,]' 24,1
$dataG[0,0] = 2.598097515653722
$dataG[1,0] = 3.755628166162433
$dataG[2,0] = 8.656069925401464
$dataG[3,0] = 8.656069925401464
$dataG[4,0] = 4.327115817150455
$dataG[5,0] = 8.656069925401464
$dataG[6,0] = 6.15379541431027
$dataG[7,0] = 1.67637130870356
$dataG[8,0] = 6.038307959104277
$dataG[9,0] = 8.656069925401464
$dataG[10,0] = 8.656069925401464
$dataG[11,0] = 4.831698963401555
$dataG[12,0] = 17.43382859256846
$dataG[13,0] = 6.15379541431027
$dataG[14,0] = 5.582307763322248
$dataG[15,0] = 16.86649396021207
$dataG[16,0] = 6.15379541431027
$dataG[17,0] = 2.998467759612273
$dataG[18,0] = 5.582307763322248
$dataG[19,0] = 4.327115817150455
$dataG[20,0] = 18.91276827552123
$dataG[21,0] = 18.91276827552123
$dataG[22,0] = 6.15379541431027
$dataG[23,0] = 8.656069925401464

$ws.Range("B2:E25").Value = $dataBE
$ws.Range("G2:G25").Value = $dataG
